$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(667, 540465, '2025-10-17T18:30:00Z', 7, 28, '1. FC Union Berlin', 18, 'Borussia Mönchengladbach', 3, 1, 'HomeWin'),
    @(668, 540462, '2025-10-18T13:30:00Z', 7, 15, '1. FSV Mainz 05', 3, 'Bayer 04 Leverkusen', 3, 4, 'AwayWin'),
    @(669, 540463, '2025-10-18T13:30:00Z', 7, 721, 'RB Leipzig', 7, 'Hamburger SV', 2, 1, 'HomeWin'),
    @(670, 540464, '2025-10-18T13:30:00Z', 7, 11, 'VfL Wolfsburg', 10, 'VfB Stuttgart', 0, 3, 'AwayWin'),
    @(671, 540467, '2025-10-18T13:30:00Z', 7, 44, '1. FC Heidenheim 1846', 12, 'SV Werder Bremen', 2, 2, 'Draw'),
    @(672, 540468, '2025-10-18T13:30:00Z', 7, 1, '1. FC Köln', 16, 'FC Augsburg', 1, 1, 'Draw'),
    @(673, 540460, '2025-10-18T16:30:00Z', 7, 5, 'FC Bayern München', 4, 'Borussia Dortmund', 2, 1, 'HomeWin'),
    @(674, 540461, '2025-10-19T13:30:00Z', 7, 17, 'SC Freiburg', 19, 'Eintracht Frankfurt', 2, 2, 'Draw'),
    @(675, 540466, '2025-10-19T15:30:00Z', 7, 20, 'FC St. Pauli 1910', 2, 'TSG 1899 Hoffenheim', 0, 3, 'AwayWin')
)


foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
}
